$wb = $excel.ActiveWorkbook

$ws2 = $wb.Worksheets.Item(2)   # "1_Vocab_Ex" - wing-shape variable classification table
$ws3 = $wb.Worksheets.Item(3)   # "2_Vocab"    - traffic-light variable classification table
$ws4 = $wb.Worksheets.Item(4)   # "3_"

# ---------------------------------------------------------------------------
# Content edits (order chosen so shared-string table is populated the same
# way the source workbook was authored)
# ---------------------------------------------------------------------------

# "2_Vocab": explanation text gets a reason clause appended
$ws3.Range("D6").Value = "This describes the state of the system at a given moment, so is a state variable."

# "1_Vocab_Ex": new explanation column (D) added
$ws2.Range("D2").Value = "This is set before the simulation is run, and remains constant."
$ws2.Range("D3").Value = "Because the simulation runs from take-off to cruising speed, this would change through time, and would describe the state of the plane at a given moment."
$ws2.Range("D6").Value = "This describes the shape of the wing, which is the variable that is being tested."

# "1_Vocab_Ex": wording corrections
$ws2.Range("A7").Value = "Speed of the aircraft just as it begins to lift off the ground"
$ws2.Range("A3").Value = "Angle of attack of the wing 10 seconds after lift-off"

$ws2.Range("D7").Value = "Because the simulation runs from take-off to cruising speed, this would change through time, and would describe the state of the plane at a given moment."

# "2_Vocab": new explanation cell added
$ws3.Range("D7").Value = "This describes the ""success"" of a particular model (presumably we'd want to minimize this!)"

# ---------------------------------------------------------------------------
# Row heights to fit the new / expanded explanatory text
# ---------------------------------------------------------------------------
$ws2.Rows.Item(2).RowHeight = 30
$ws2.Rows.Item(3).RowHeight = 75
$ws2.Rows.Item(6).RowHeight = 45
$ws2.Rows.Item(7).RowHeight = 75

$ws3.Rows.Item(6).RowHeight = 45
$ws3.Rows.Item(7).RowHeight = 45

# ---------------------------------------------------------------------------
# Selection / cursor position ($ws3 stays active/selected last so the
# workbook's active-tab pointer keeps pointing at "2_Vocab", matching the
# source which never changed the active sheet)
# ---------------------------------------------------------------------------
$ws2.Range("A4").Select()
$ws4.Range("A1").Select()
$ws3.Range("D8").Select()
